# Generate Report for Handoff
# New handoff/report-generation timestamps are recorded for the last tracked
# file (b6bcf457-6b32-4650-b0e2-0b6de2e94e97) across the per-locale sheets,
# and the Overview roll-up column is refreshed to the newest of those dates.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# "Latest Handoff Datetime" column (H) for the last data row (row 7) of each
# locale sheet.
$zhcn.Range("H7").Value = "2016-09-05 00:47:17"
$dede.Range("H7").Value = "2016-09-05 00:47:21"

# "Latest HO Xliff Generate Date" column (G) on the Overview sheet mirrors the
# newest per-locale handoff datetime for that row.
$overview.Range("G7").Value = "2016-09-05 00:47:21"
